# daily auto push: 2025-10-06 02:01 UTC
# Appends the new daily log entry as row 69 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date-like text (e.g. "2025/10/06") that must remain a
# literal string instead of being auto-converted into a date serial value.
# Force the cell to text format before assigning, then clear the
# formatting afterwards so the new row matches the plain (unstyled) look
# of the rest of the data rows.
$ws.Range("A69").NumberFormat = "@"
$ws.Range("A69").Value = "2025/10/06"
$ws.Range("A69").ClearFormats()

$ws.Range("B69").Value = "月"
$ws.Range("C69").Value = 8
$ws.Range("D69").Value = 201
